# Automatische test-sync: 2025-06-29 14:04:50
# Adds a new test-mail row to the "Logs" sheet, rolls the matching
# "Productinformatie" category count into the "Dashboard" sheet, and
# extends the dashboard bar chart + conditional formatting ranges to
# cover the newly added row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 6 with the new test mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(6, 1).Value = "Wat zijn de verzendkosten?"
$logs.Cells.Item(6, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(6, 3).Value = "Testmail #3: Wat zijn de verzendkosten?"
$logs.Cells.Item(6, 4).Value = "Productinformatie"
$logs.Cells.Item(6, 5).Value = "Beste afzender,`nDank u voor uw interesse in onze producten/diensten. Om u nauwkeurige informatie te verstrekken over de verzendkosten, hebben we meer details nodig zoals het product/dienst waar u naar informeert en het afleveradres. Zou u ons kunnen voorzien van deze informatie zodat we u een precieze schatting van de verzendkosten kunnen geven?`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Cells.Item(6, 6).Value = "2025-06-29 14:04:02"
$logs.Cells.Item(6, 7).Value = "Ja"
$logs.Cells.Item(6, 8).Value = "Nee"
$logs.Cells.Item(6, 9).Value = "Ja"

# Extend the conditional-formatting ranges (D/G/H/I 2:5 -> 2:6) so the new
# row is covered, without disturbing the rule definitions themselves.
$logs.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))
$logs.Range("H2:H5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H6"))
$logs.Range("I2:I5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I6"))

# ---------------------------------------------------------------------
# 2) Dashboard sheet: append the "Productinformatie" roll-up row
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Productinformatie"
$dash.Cells.Item(4, 2).Value = 1

# ---------------------------------------------------------------------
# 3) Dashboard chart: widen the category/value series ranges to include
#    the newly added row 4
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$4,'Dashboard'!`$B`$2:`$B`$4,1)"
